# Auto-generated Excel COM-interop script
# Applies a scheduled data-refresh update to the Leve profit sheets
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 48
$ws.Cells.Item(5, 9).Value = 48
$ws.Cells.Item(5, 11).Value = 48
$ws.Cells.Item(5, 13).Value = 67
$ws.Cells.Item(12, 8).Value = 225.81818
$ws.Cells.Item(12, 9).Value = 220.44444
$ws.Cells.Item(12, 11).Value = 220.44444
$ws.Cells.Item(12, 13).Value = -50.44443999999999
$ws.Cells.Item(17, 8).Value = 1344.619
$ws.Cells.Item(17, 10).Value = 1351.5278
$ws.Cells.Item(17, 12).Value = 4054.5834
$ws.Cells.Item(17, 14).Value = -4390.5834
$ws.Cells.Item(48, 8).Value = 4974.5913
$ws.Cells.Item(48, 10).Value = 4974.5913
$ws.Cells.Item(48, 12).Value = 14923.7739
$ws.Cells.Item(48, 14).Value = -15507.7739
$ws.Cells.Item(56, 8).Value = 4974.5913
$ws.Cells.Item(56, 10).Value = 4974.5913
$ws.Cells.Item(56, 12).Value = 14923.7739
$ws.Cells.Item(56, 14).Value = -15991.7739

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4703.524
$ws.Cells.Item(2, 9).Value = 5898.364
$ws.Cells.Item(2, 11).Value = 5898.364
$ws.Cells.Item(2, 13).Value = -5785.364
$ws.Cells.Item(5, 8).Value = 4339.4287
$ws.Cells.Item(5, 9).Value = 4594.5
$ws.Cells.Item(5, 11).Value = 4594.5
$ws.Cells.Item(5, 13).Value = -4482.5
$ws.Cells.Item(32, 8).Value = 7648.106
$ws.Cells.Item(32, 9).Value = 7611.923
$ws.Cells.Item(32, 10).Value = 10000
$ws.Cells.Item(32, 11).Value = 7611.923
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = -7324.923
$ws.Cells.Item(32, 14).Value = -10574
$ws.Cells.Item(45, 8).Value = 5582.52
$ws.Cells.Item(45, 9).Value = 4883.952
$ws.Cells.Item(45, 11).Value = 4883.952
$ws.Cells.Item(45, 13).Value = -4506.952
$ws.Cells.Item(61, 8).Value = 3913.1538
$ws.Cells.Item(61, 9).Value = 3629.682
$ws.Cells.Item(61, 10).Value = 4507.095
$ws.Cells.Item(61, 11).Value = 3629.682
$ws.Cells.Item(61, 12).Value = 4507.095
$ws.Cells.Item(61, 13).Value = -3417.682
$ws.Cells.Item(61, 14).Value = -4931.095
$ws.Cells.Item(74, 8).Value = 2685.6484
$ws.Cells.Item(74, 9).Value = 3825.4614
$ws.Cells.Item(74, 10).Value = 1830.7885
$ws.Cells.Item(74, 11).Value = 3825.4614
$ws.Cells.Item(74, 12).Value = 1830.7885
$ws.Cells.Item(74, 13).Value = -2951.4614
$ws.Cells.Item(74, 14).Value = -3578.7885
$ws.Cells.Item(77, 8).Value = 2685.6484
$ws.Cells.Item(77, 9).Value = 3825.4614
$ws.Cells.Item(77, 10).Value = 1830.7885
$ws.Cells.Item(77, 11).Value = 19127.307
$ws.Cells.Item(77, 12).Value = 9153.942500000001
$ws.Cells.Item(77, 13).Value = -14759.307
$ws.Cells.Item(77, 14).Value = -17889.9425
$ws.Cells.Item(88, 8).Value = 83335680
$ws.Cells.Item(88, 9).Value = 1933
$ws.Cells.Item(88, 10).Value = 100002424
$ws.Cells.Item(88, 11).Value = 1933
$ws.Cells.Item(88, 12).Value = 100002424
$ws.Cells.Item(88, 13).Value = -1527
$ws.Cells.Item(88, 14).Value = -100003236
$ws.Cells.Item(91, 8).Value = 83335680
$ws.Cells.Item(91, 9).Value = 1933
$ws.Cells.Item(91, 10).Value = 100002424
$ws.Cells.Item(91, 11).Value = 1933
$ws.Cells.Item(91, 12).Value = 100002424
$ws.Cells.Item(91, 13).Value = -529
$ws.Cells.Item(91, 14).Value = -100005232
$ws.Cells.Item(97, 8).Value = 20013008
$ws.Cells.Item(97, 9).Value = 19364
$ws.Cells.Item(97, 11).Value = 19364
$ws.Cells.Item(97, 13).Value = -18868
$ws.Cells.Item(116, 8).Value = 4703.524
$ws.Cells.Item(116, 9).Value = 5898.364
$ws.Cells.Item(116, 11).Value = 5898.364
$ws.Cells.Item(116, 13).Value = -3604.364
$ws.Cells.Item(122, 8).Value = 861589.1
$ws.Cells.Item(122, 9).Value = 4021.4138
$ws.Cells.Item(122, 10).Value = 5006499.5
$ws.Cells.Item(122, 11).Value = 12064.2414
$ws.Cells.Item(122, 12).Value = 15019498.5
$ws.Cells.Item(122, 13).Value = -9614.241399999999
$ws.Cells.Item(122, 14).Value = -15024398.5
$ws.Cells.Item(136, 8).Value = 3913.1538
$ws.Cells.Item(136, 9).Value = 3629.682
$ws.Cells.Item(136, 10).Value = 4507.095
$ws.Cells.Item(136, 11).Value = 10889.046
$ws.Cells.Item(136, 12).Value = 13521.285
$ws.Cells.Item(136, 13).Value = -8339.045999999998
$ws.Cells.Item(136, 14).Value = -18621.285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4703.524
$ws.Cells.Item(3, 9).Value = 5898.364
$ws.Cells.Item(3, 11).Value = 5898.364
$ws.Cells.Item(3, 13).Value = -5784.364
$ws.Cells.Item(4, 8).Value = 4339.4287
$ws.Cells.Item(4, 9).Value = 4594.5
$ws.Cells.Item(4, 11).Value = 4594.5
$ws.Cells.Item(4, 13).Value = -4479.5
$ws.Cells.Item(105, 8).Value = 12064.286
$ws.Cells.Item(105, 9).Value = 12990
$ws.Cells.Item(105, 10).Value = 9750
$ws.Cells.Item(105, 11).Value = 12990
$ws.Cells.Item(105, 12).Value = 9750
$ws.Cells.Item(105, 13).Value = -11243
$ws.Cells.Item(105, 14).Value = -13244
$ws.Cells.Item(107, 8).Value = 1913.2
$ws.Cells.Item(107, 9).Value = 1807.7059
$ws.Cells.Item(107, 11).Value = 1807.7059
$ws.Cells.Item(107, 13).Value = 112.2941000000001
$ws.Cells.Item(134, 8).Value = 1986.0834
$ws.Cells.Item(134, 9).Value = 1279.8948
$ws.Cells.Item(134, 11).Value = 3839.6844
$ws.Cells.Item(134, 13).Value = -1304.6844

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 329.9375
$ws.Cells.Item(7, 9).Value = 244.53847
$ws.Cells.Item(7, 11).Value = 244.53847
$ws.Cells.Item(7, 13).Value = -131.53847
$ws.Cells.Item(22, 8).Value = 1017.38464
$ws.Cells.Item(22, 9).Value = 716.375
$ws.Cells.Item(22, 10).Value = 1499
$ws.Cells.Item(22, 11).Value = 716.375
$ws.Cells.Item(22, 12).Value = 1499
$ws.Cells.Item(22, 13).Value = -366.375
$ws.Cells.Item(22, 14).Value = -2199
$ws.Cells.Item(58, 8).Value = 3028.85
$ws.Cells.Item(58, 10).Value = 6139
$ws.Cells.Item(58, 12).Value = 6139
$ws.Cells.Item(58, 14).Value = -6545
$ws.Cells.Item(105, 8).Value = 11330.637
$ws.Cells.Item(105, 9).Value = 12323.7
$ws.Cells.Item(105, 10).Value = 1400
$ws.Cells.Item(105, 11).Value = 12323.7
$ws.Cells.Item(105, 12).Value = 1400
$ws.Cells.Item(105, 13).Value = -10576.7
$ws.Cells.Item(105, 14).Value = -4894
$ws.Cells.Item(132, 8).Value = 12753.444
$ws.Cells.Item(132, 9).Value = 17833
$ws.Cells.Item(132, 11).Value = 53499
$ws.Cells.Item(132, 13).Value = -50969
$ws.Cells.Item(136, 8).Value = 3028.85
$ws.Cells.Item(136, 10).Value = 6139
$ws.Cells.Item(136, 12).Value = 18417
$ws.Cells.Item(136, 14).Value = -23517
$ws.Cells.Item(141, 8).Value = 260818.81
$ws.Cells.Item(141, 10).Value = 309756.25
$ws.Cells.Item(141, 12).Value = 309756.25
$ws.Cells.Item(141, 14).Value = -320116.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 322.9
$ws.Cells.Item(2, 9).Value = 304.2857
$ws.Cells.Item(2, 11).Value = 1825.7142
$ws.Cells.Item(2, 13).Value = -1712.7142
$ws.Cells.Item(5, 8).Value = 590901.6
$ws.Cells.Item(5, 10).Value = 1115221.5
$ws.Cells.Item(5, 12).Value = 3345664.5
$ws.Cells.Item(5, 14).Value = -3345888.5
$ws.Cells.Item(122, 8).Value = 2206.2
$ws.Cells.Item(122, 9).Value = 1139.5333
$ws.Cells.Item(122, 10).Value = 2846.2
$ws.Cells.Item(122, 11).Value = 10255.7997
$ws.Cells.Item(122, 12).Value = 25615.8
$ws.Cells.Item(122, 13).Value = -7805.7997
$ws.Cells.Item(122, 14).Value = -30515.8
$ws.Cells.Item(125, 8).Value = 30
$ws.Cells.Item(125, 9).Value = 30
$ws.Cells.Item(125, 11).Value = 90
$ws.Cells.Item(125, 13).Value = 4830
$ws.Cells.Item(135, 8).Value = 590901.6
$ws.Cells.Item(135, 10).Value = 1115221.5
$ws.Cells.Item(135, 12).Value = 10036993.5
$ws.Cells.Item(135, 14).Value = -10042063.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 8388.954
$ws.Cells.Item(102, 9).Value = 9392.056
$ws.Cells.Item(102, 10).Value = 3875
$ws.Cells.Item(102, 11).Value = 9392.056
$ws.Cells.Item(102, 12).Value = 3875
$ws.Cells.Item(102, 13).Value = -7770.056
$ws.Cells.Item(102, 14).Value = -7119
$ws.Cells.Item(132, 8).Value = 3017.5186
$ws.Cells.Item(132, 9).Value = 3011.375
$ws.Cells.Item(132, 10).Value = 3035.0715
$ws.Cells.Item(132, 11).Value = 9034.125
$ws.Cells.Item(132, 12).Value = 9105.2145
$ws.Cells.Item(132, 13).Value = -6504.125
$ws.Cells.Item(132, 14).Value = -14165.2145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 7296.864
$ws.Cells.Item(16, 9).Value = 8089.278
$ws.Cells.Item(16, 10).Value = 3731
$ws.Cells.Item(16, 11).Value = 8089.278
$ws.Cells.Item(16, 12).Value = 3731
$ws.Cells.Item(16, 13).Value = -7919.278
$ws.Cells.Item(16, 14).Value = -4071
$ws.Cells.Item(25, 8).Value = 9666.333000000001
$ws.Cells.Item(25, 9).Value = 8999
$ws.Cells.Item(25, 11).Value = 8999
$ws.Cells.Item(25, 13).Value = -8769
$ws.Cells.Item(40, 8).Value = 22130.94
$ws.Cells.Item(40, 9).Value = 27754.264
$ws.Cells.Item(40, 11).Value = 27754.264
$ws.Cells.Item(40, 13).Value = -27618.264
$ws.Cells.Item(46, 8).Value = 1411.2222
$ws.Cells.Item(46, 10).Value = 2633.6667
$ws.Cells.Item(46, 12).Value = 2633.6667
$ws.Cells.Item(46, 14).Value = -3009.6667
$ws.Cells.Item(68, 8).Value = 5334
$ws.Cells.Item(68, 9).Value = 3037.375
$ws.Cells.Item(68, 10).Value = 6747.3076
$ws.Cells.Item(68, 11).Value = 3037.375
$ws.Cells.Item(68, 12).Value = 6747.3076
$ws.Cells.Item(68, 13).Value = -2288.375
$ws.Cells.Item(68, 14).Value = -8245.3076
$ws.Cells.Item(71, 8).Value = 5334
$ws.Cells.Item(71, 9).Value = 3037.375
$ws.Cells.Item(71, 10).Value = 6747.3076
$ws.Cells.Item(71, 11).Value = 15186.875
$ws.Cells.Item(71, 12).Value = 33736.538
$ws.Cells.Item(71, 13).Value = -11442.875
$ws.Cells.Item(71, 14).Value = -41224.538
$ws.Cells.Item(93, 8).Value = 3594.889
$ws.Cells.Item(93, 9).Value = 3759.7273
$ws.Cells.Item(93, 10).Value = 2869.6
$ws.Cells.Item(93, 11).Value = 3759.7273
$ws.Cells.Item(93, 12).Value = 2869.6
$ws.Cells.Item(93, 13).Value = -2511.7273
$ws.Cells.Item(93, 14).Value = -5365.6
$ws.Cells.Item(100, 8).Value = 4933.8184
$ws.Cells.Item(100, 9).Value = 1878.8334
$ws.Cells.Item(100, 11).Value = 1878.8334
$ws.Cells.Item(100, 13).Value = -1337.8334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 165500000
$ws.Cells.Item(94, 9).Value = 1000000
$ws.Cells.Item(94, 11).Value = 1000000
$ws.Cells.Item(94, 13).Value = -999099
$ws.Cells.Item(122, 8).Value = 4275.7427
$ws.Cells.Item(122, 9).Value = 2354.4348
$ws.Cells.Item(122, 11).Value = 7063.3044
$ws.Cells.Item(122, 13).Value = -4613.3044
$ws.Cells.Item(132, 8).Value = 4940.84
$ws.Cells.Item(132, 9).Value = 5350.515
$ws.Cells.Item(132, 10).Value = 4145.5884
$ws.Cells.Item(132, 11).Value = 16051.545
$ws.Cells.Item(132, 12).Value = 12436.7652
$ws.Cells.Item(132, 13).Value = -13521.545
$ws.Cells.Item(132, 14).Value = -17496.7652
$ws.Cells.Item(136, 8).Value = 419078.7
$ws.Cells.Item(136, 9).Value = 469407.9
$ws.Cells.Item(136, 11).Value = 1408223.7
$ws.Cells.Item(136, 13).Value = -1405673.7

Write-Host "Applied $([int]253) cell updates across 8 sheets"
